$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "56.944.00"
$ws.Cells.Item(2, 5).Value = "  -0.65%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.314.85"
$ws.Cells.Item(3, 5).Value = "  -2.11%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'529.16"
$ws.Cells.Item(5, 5).Value = "  +1.64%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'132.37"
$ws.Cells.Item(6, 5).Value = "  -2.61%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.994"
$ws.Cells.Item(7, 5).Value = "  -0.19%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -1.31%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.337.28"
$ws.Cells.Item(9, 5).Value = "  -2.07%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -1.69%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.02%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'5.34"

# Row 13
$ws.Cells.Item(13, 4).Value = "'0.351"
$ws.Cells.Item(13, 5).Value = "  +2.21%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.733.89"
$ws.Cells.Item(14, 5).Value = "  -1.83%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'23.35"
$ws.Cells.Item(15, 5).Value = "  -4.54%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "56.951.49"
$ws.Cells.Item(16, 5).Value = "  -0.65%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -2.26%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.339.00"
$ws.Cells.Item(18, 5).Value = "  -1.60%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'335.22"
$ws.Cells.Item(19, 5).Value = "  +1.55%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'10.41"
$ws.Cells.Item(20, 5).Value = "  -1.68%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -1.86%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'6.80"
$ws.Cells.Item(22, 5).Value = "  +1.21%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.44%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'61.90"
$ws.Cells.Item(24, 5).Value = "  +0.69%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'8.69"
$ws.Cells.Item(26, 5).Value = "  -4.16%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'0.992"
$ws.Cells.Item(27, 5).Value = "  -2.26%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.64%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'173.27"
$ws.Cells.Item(29, 5).Value = "  +3.88%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'1.72"
$ws.Cells.Item(30, 5).Value = "  +0.14%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.0₃0723"
$ws.Cells.Item(31, 5).Value = "  -3.52%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'6.10"
$ws.Cells.Item(32, 5).Value = "  -3.05%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'18.46"
$ws.Cells.Item(33, 5).Value = "  -0.86%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'0.998"

# Row 35
$ws.Cells.Item(35, 5).Value = "  -0.26%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -4.14%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.59%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -1.70%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "OKB"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(39, 4).Value = "'39.14"
$ws.Cells.Item(39, 5).Value = "  +0.95%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(40, 4).Value = "'1.58"
$ws.Cells.Item(40, 5).Value = "  -1.86%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "RenderToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(41, 4).Value = "'5.77"
$ws.Cells.Item(41, 5).Value = "  +8.06%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'148.62"
$ws.Cells.Item(42, 5).Value = "  -0.87%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -3.07%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -1.95%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'283.03"
$ws.Cells.Item(45, 5).Value = "  -3.18%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.0932"
$ws.Cells.Item(46, 5).Value = "  -0.88%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -2.00%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'18.73"
$ws.Cells.Item(48, 5).Value = "  +2.80%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -1.77%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "VeChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(50, 4).Value = "'0.0216"
$ws.Cells.Item(50, 5).Value = "  -1.51%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Polygon"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(51, 4).Value = "'0.382"
$ws.Cells.Item(51, 5).Value = "  +5.79%  "
